# Commit: "review an IELTS essay and add new words"
# Add a new "cloth / clothe" mispronunciation entry to the word list on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Enter the two new words first (cloth, clothe), then their pronunciations,
# matching the order the shared-string table was built in.
$ws.Range("A6").Value = "cloth"
$ws.Range("A7").Value = "clothe"

$ws.Range("B6").Value = "/KlOT/ n."
$ws.Range("C6").Value = "/klEuT/"

$ws.Range("B7").Value = "/klED/ vt."

# The "mispronounce" column for "cloth" (C6) is shown struck-through,
# matching the formatting used for the other mispronunciation entries.
$ws.Range("C6").Font.Strikethrough = $true

# Make sure Sheet2 stays the active sheet and update the selected cell.
$ws.Activate() | Out-Null
$ws.Range("B7").Select() | Out-Null
